$wb = $excel.ActiveWorkbook

# --- Fix the Germany sheet's stale "select all" selection (A1:XFD1048576 -> A1:D11) ---
$germany = $wb.Worksheets.Item("Germany")
$germany.Range("A1:D11").Select()

# --- Duplicate the Swiss sheet to create the new "Portugal" sheet, placed right after Swiss ---
$swiss = $wb.Worksheets.Item("Swiss")
$swiss.Copy($null, $swiss)
$portugal = $wb.Worksheets.Item($wb.Worksheets.Count)
$portugal.Name = "Portugal"

# --- Update the Portugal-specific content ---
$portugal.Range("B2").Value = "Portugal Market"
$portugal.Range("B4").Value = "NGC-3479/T2459/T-2460/T2461/T2462"

# --- Column widths specific to the Portugal sheet ---
$portugal.Columns.Item(1).ColumnWidth = 26.944010416666668
$portugal.Columns.Item(2).ColumnWidth = 32.276041666666664
$portugal.Columns.Item(3).ColumnWidth = 12.721354166666666
$portugal.Columns.Item(4).ColumnWidth = 11.166666666666666

# --- Rows 3, 4 & 5 grow taller on the Portugal sheet ---
$portugal.Rows.Item(3).RowHeight = 28.8
$portugal.Rows.Item(4).RowHeight = 28.8
$portugal.Rows.Item(5).RowHeight = 28.8

# --- Selection/active cell on Portugal is B4 ---
$portugal.Range("B4").Select()

# --- Make Portugal the active sheet/tab ---
$portugal.Activate()
